$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (shared-string text updates)
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Update GDP (column C) values with higher-precision figures
$ws.Range("C2").Value = 5596.139681459835
$ws.Range("C3").Value = 1873.394108966653
$ws.Range("C4").Value = 1909.084588129339
$ws.Range("C5").Value = 6128.19547247793
$ws.Range("C6").Value = 4729.735976516416
$ws.Range("C7").Value = 1268.249210347625
$ws.Range("C8").Value = 1286.515571617672
$ws.Range("C9").Value = 513.4456986202239
$ws.Range("C10").Value = 401.8350013668368
$ws.Range("C11").Value = 951.6879611168786
$ws.Range("C12").Value = 5730.354774594881
$ws.Range("C13").Value = 1904.346464968814
$ws.Range("C14").Value = 9271.398233246389
$ws.Range("C15").Value = 1955.461557360978
$ws.Range("C16").Value = 4633.590358399045
$ws.Range("C17").Value = 5082.354756663512
$ws.Range("C18").Value = 2217.474008566157
$ws.Range("C19").Value = 1357.563719132622
$ws.Range("C20").Value = 1037.747039954749
$ws.Range("C21").Value = 1446.371630707023
$ws.Range("C22").Value = 538.1162350013198
$ws.Range("C23").Value = 1263.452411343738
$ws.Range("C24").Value = 417.6031683854853
$ws.Range("C25").Value = 553.2014555484933
$ws.Range("C26").Value = 5885.254624554112
$ws.Range("C27").Value = 1939.33862702996
$ws.Range("C28").Value = 9477.887185090232
$ws.Range("C29").Value = 10883.31535948899
$ws.Range("C30").Value = 2024.117324382548
$ws.Range("C31").Value = 6711.616186806423
$ws.Range("C32").Value = 4921.848409120176
$ws.Range("C33").Value = 5360.226632400601
$ws.Range("C34").Value = 4961.234688573883
$ws.Range("C35").Value = 573.9239887389259
$ws.Range("C36").Value = 1325.930225429421
$ws.Range("C37").Value = 564.5208442217756
$ws.Range("C38").Value = 1543.763984230257
$ws.Range("C39").Value = 1291.622214254295
$ws.Range("C40").Value = 1469.177610078392
$ws.Range("C41").Value = 470.3014046213344
$ws.Range("C42").Value = 5122.180090208862
$ws.Range("C43").Value = 5642.578115155247
$ws.Range("C44").Value = 2094.024217383061
$ws.Range("C45").Value = 1982.009737844954
$ws.Range("C46").Value = 6051.685746144485
$ws.Range("C47").Value = 593.1620921048029
$ws.Range("C48").Value = 1360.10887014004
$ws.Range("C49").Value = 597.3813896804552
$ws.Range("C50").Value = 1618.597849849475
$ws.Range("C51").Value = 1291.415042301529
$ws.Range("C52").Value = 1544.619247249133
$ws.Range("C53").Value = 555.2055623950326
$ws.Range("C54").Value = 5295.682695961288
$ws.Range("C55").Value = 5919.20956823756
$ws.Range("C56").Value = 2201.396847776877
$ws.Range("C57").Value = 2000.792448761861
$ws.Range("C58").Value = 6203.843262938323
$ws.Range("C59").Value = 10398.69400694643
$ws.Range("C60").Value = 2286.013198234259
$ws.Range("C61").Value = 1401.753174264641
$ws.Range("C62").Value = 7449.08671983612
$ws.Range("C63").Value = 6255.426161047989
$ws.Range("C64").Value = 2612.856880840196
$ws.Range("C65").Value = 1627.760281433693
$ws.Range("C66").Value = 1640.18070024053
$ws.Range("C67").Value = 586.2293607842975
$ws.Range("C68").Value = 558.2093442539386
$ws.Range("C69").Value = 5996.49696468919
$ws.Range("C70").Value = 1338.716747746975
$ws.Range("C71").Value = 567.8342670439314
$ws.Range("C72").Value = 10568.15780870825
$ws.Range("C73").Value = 2361.056581219794
$ws.Range("C74").Value = 1441.783971398429
$ws.Range("C75").Value = 7580.275568826287
$ws.Range("C76").Value = 6522.736799041846
$ws.Range("C77").Value = 2735.187532014817
$ws.Range("C78").Value = 1625.905825842452
$ws.Range("C79").Value = 1751.664428859304
$ws.Range("C80").Value = 571.453129531788
$ws.Range("C81").Value = 579.0880693780265
$ws.Range("C82").Value = 6114.227214287786
$ws.Range("C83").Value = 1384.519227335143
$ws.Range("C84").Value = 441.1376640642927
$ws.Range("C85").Value = 10239.48134799327
$ws.Range("C86").Value = 2425.561644739583
$ws.Range("C87").Value = 1469.192636109792
$ws.Range("C88").Value = 7633.969039669125
$ws.Range("C89").Value = 6550.274372976741
$ws.Range("C90").Value = 5176.058803160127
$ws.Range("C91").Value = 1644.598009122967
$ws.Range("C92").Value = 2111.193164269742
$ws.Range("C93").Value = 1875.732161108182
$ws.Range("C94").Value = 548.2681436079887
$ws.Range("C95").Value = 6262.368904654469
$ws.Range("C96").Value = 1431.756130822538
$ws.Range("C97").Value = 457.8330917196623
$ws.Range("C98").Value = 0
$ws.Range("C99").Value = 2448.861248735403
